$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This engine's Hyperlinks collection only supports clearing every
# hyperlink on the sheet at once (deleting a single item removes them
# all), so drop all three, delete the old rows, then recreate the one
# hyperlink that should survive (A2 -> appumuv@gmail.com).
$ws.Hyperlinks.Delete()

# Drop the old "abc@abc.com" / "123@113.com" test rows entirely.
$ws.Rows("3:4").Delete()

$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:appumuv@gmail.com")
$ws.Range("A2").Style = "Hyperlink"

# Rework the remaining two rows into the new "Product / Quantity / location"
# end-to-end test data.
$ws.Range("C1").Value = "Product"
$ws.Range("D1").Value = "Quantity"
$ws.Range("E1").Value = "location"

$ws.Range("E2").Value = "Cognizant"
$ws.Range("C2").Value = "Noise Grey Knitted Slouchy Beanie"
$ws.Range("D2").Value = 2

$ws.Columns("C:C").AutoFit()

Write-Output "done"
